$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StepperPage")

# Update the username/email cell and its dependent look-alike hyperlink styling.
$ws.Range("B2").Value = "automationtest@pixentia.com"

# Update the two "test" placeholder cells with the new Automation/Test values.
$ws.Range("D2").Value = "Automation"
$ws.Range("E2").Value = "Test"

# Add the mailto hyperlink on B2 and apply the built-in Hyperlink style.
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:automationtest@pixentia.com")
$ws.Range("B2").Style = "Hyperlink"

# Update the active selection shown when the sheet is opened.
$ws.Range("E2").Select()
